# Weekly update: a new price record was reported for "Haba" (Femacal de
# La Calera) and inserted as a new row right after the current row 70,
# pushing all subsequent rows (old 71..91) down by one (new 72..92).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 71; existing rows 71-91 shift to 72-92.
$ws.Rows("71:71").Insert()

# Populate the newly inserted row 71 with the new record.
$ws.Range("A71").Value = 3
$ws.Range("B71").Value = "Femacal de La Calera"
$ws.Range("C71").Value = "Coquimbo"
$ws.Range("D71").Value = 44511
$ws.Range("E71").Value = 5
$ws.Range("F71").Value = 100112026
$ws.Range("G71").Value = "Haba"
$ws.Range("H71").Value = "Sin especificar"
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 80
$ws.Range("K71").Value = 7000
$ws.Range("L71").Value = 7500
$ws.Range("M71").Value = 7250
$ws.Range("N71").Value = "`$/saco 25 kilos"
$ws.Range("O71").Value = "Provincia de Quillota"
$ws.Range("P71").Value = 290
$ws.Range("Q71").Value = 25
$ws.Range("R71").Value = "Hortaliza"

# Make sure the date cell keeps the same date/time number format used by
# the rest of the "Fecha" column.
$ws.Range("D71").NumberFormat = "YYYY-MM-DD HH:MM:SS"
